$p = $ppt.ActivePresentation

function Find-ShapeWithText {
    param($slide, [string]$needle)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -and $t.Contains($needle)) {
                return $shp
            }
        }
    }
    return $null
}

function Replace-RunText {
    param($textRange, [string]$oldText, [string]$newText)
    $full = $textRange.Text
    $idx0 = $full.IndexOf($oldText)
    if ($idx0 -lt 0) {
        return $false
    }
    $start1 = $idx0 + 1
    $sub = $textRange.Characters($start1, $oldText.Length)
    $sub.Text = $newText
    return $true
}

$zwsp = [char]0x200b
$endash = [char]0x2013

# --- Slide "Early Approaches:" / "Deep Learning Advances:" ---------------
# Hybrid Models (e.g., AlexNet + Grey Wolf Optimizer) accuracy figure
# updated from >99% to >90%.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $shp = Find-ShapeWithText $slide "Grey Wolf Optimizer"
    if ($shp -ne $null) {
        $tr = $shp.TextFrame.TextRange
        $oldRun = " + Grey Wolf Optimizer) achieved >99% accuracy on ISIC datasets" + $zwsp + "."
        $newRun = " + Grey Wolf Optimizer) achieved >90% accuracy on ISIC datasets" + $zwsp + "."
        Replace-RunText $tr $oldRun $newRun | Out-Null
    }
}

# --- Slide "Emerging Trends:" ---------------------------------------------
# Drop the "(up to 96%-100%)" qualifier from the ViT accuracy bullet.
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    $shp = Find-ShapeWithText $slide "HAM10000"
    if ($shp -ne $null) {
        $tr = $shp.TextFrame.TextRange
        $oldRun = "Achieved high accuracy (up to 96%" + $endash + "100%) on datasets like HAM10000" + $zwsp + "."
        $newRun = "Achieved high accuracy on datasets like HAM10000" + $zwsp + "."
        Replace-RunText $tr $oldRun $newRun | Out-Null
    }
}
